# Update "想去人数" (interest count) values in the 展览 (Exhibition), 演出
# (Performance), and 全部类型 (All types) sheets to reflect newly generated
# stats, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibition.Range("F13").Value = 6586
$wsExhibition.Range("F14").Value = 41
$wsExhibition.Range("F15").Value = 88
$wsExhibition.Range("F17").Value = 5043
$wsExhibition.Range("F20").Value = 4215
$wsExhibition.Range("F22").Value = 4155
$wsExhibition.Range("F33").Value = 7529
$wsExhibition.Range("F39").Value = 65
$wsExhibition.Range("F40").Value = 1502
$wsExhibition.Range("F42").Value = 840
$wsExhibition.Range("F44").Value = 3696
$wsExhibition.Range("F49").Value = 1042

# --- 演出 (sheet2) ---
$wsPerformance.Range("F21").Value = 857

# --- 全部类型 (sheet4) ---
$wsAll.Range("F16").Value = 6586
$wsAll.Range("F17").Value = 41
$wsAll.Range("F18").Value = 88
$wsAll.Range("F20").Value = 5043
$wsAll.Range("F23").Value = 4215
$wsAll.Range("F24").Value = 4155
$wsAll.Range("F33").Value = 7529
$wsAll.Range("F39").Value = 65
$wsAll.Range("F40").Value = 1502
$wsAll.Range("F42").Value = 840
$wsAll.Range("F44").Value = 3696
$wsAll.Range("F48").Value = 1042
